# Auto-generated Excel COM-interop script to apply profit/price recalculations
# to the Zeromus_Profits workbook, one sheet worksheet at a time.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(6, 8).Value = 841.875
$ws.Cells.Item(6, 9).Value = 968.3333
$ws.Cells.Item(6, 10).Value = 462.5
$ws.Cells.Item(6, 11).Value = 2904.9999
$ws.Cells.Item(6, 12).Value = 1387.5
$ws.Cells.Item(6, 13).Value = -2792.9999
$ws.Cells.Item(6, 14).Value = -1611.5
$ws.Cells.Item(53, 8).Value = 426.7857
$ws.Cells.Item(53, 9).Value = 569.25
$ws.Cells.Item(53, 10).Value = 236.83333
$ws.Cells.Item(53, 11).Value = 569.25
$ws.Cells.Item(53, 12).Value = 236.83333
$ws.Cells.Item(53, 13).Value = 67.75
$ws.Cells.Item(53, 14).Value = -1510.83333
$ws.Cells.Item(74, 8).Value = 4609.727
$ws.Cells.Item(74, 9).Value = 4958.143
$ws.Cells.Item(74, 10).Value = 4000
$ws.Cells.Item(74, 11).Value = 4958.143
$ws.Cells.Item(74, 12).Value = 4000
$ws.Cells.Item(74, 13).Value = -4022.143
$ws.Cells.Item(74, 14).Value = -5872
$ws.Cells.Item(77, 8).Value = 4609.727
$ws.Cells.Item(77, 9).Value = 4958.143
$ws.Cells.Item(77, 10).Value = 4000
$ws.Cells.Item(77, 11).Value = 24790.715
$ws.Cells.Item(77, 12).Value = 20000
$ws.Cells.Item(77, 13).Value = -20110.715
$ws.Cells.Item(77, 14).Value = -29360
$ws.Cells.Item(88, 8).Value = 7374.4443
$ws.Cells.Item(88, 9).Value = 21450
$ws.Cells.Item(88, 10).Value = 3352.8572
$ws.Cells.Item(88, 11).Value = 21450
$ws.Cells.Item(88, 12).Value = 3352.8572
$ws.Cells.Item(88, 13).Value = -21044
$ws.Cells.Item(88, 14).Value = -4164.8572
$ws.Cells.Item(91, 8).Value = 7374.4443
$ws.Cells.Item(91, 9).Value = 21450
$ws.Cells.Item(91, 10).Value = 3352.8572
$ws.Cells.Item(91, 11).Value = 21450
$ws.Cells.Item(91, 12).Value = 3352.8572
$ws.Cells.Item(91, 13).Value = -20046
$ws.Cells.Item(91, 14).Value = -6160.8572
$ws.Cells.Item(111, 8).Value = 2650.889
$ws.Cells.Item(111, 9).Value = 2465.4285
$ws.Cells.Item(111, 10).Value = 3300
$ws.Cells.Item(111, 11).Value = 7396.2855
$ws.Cells.Item(111, 12).Value = 9900
$ws.Cells.Item(111, 13).Value = -4329.2855
$ws.Cells.Item(111, 14).Value = -16034
$ws.Cells.Item(112, 8).Value = 1147
$ws.Cells.Item(112, 10).Value = 1307.0714
$ws.Cells.Item(112, 12).Value = 3921.2142
$ws.Cells.Item(112, 14).Value = -6137.2142
$ws.Cells.Item(127, 8).Value = 1157.1428
$ws.Cells.Item(127, 9).Value = 425
$ws.Cells.Item(127, 11).Value = 1275
$ws.Cells.Item(127, 13).Value = 3685
$ws.Cells.Item(129, 8).Value = 1108.28
$ws.Cells.Item(129, 10).Value = 1301.15
$ws.Cells.Item(129, 12).Value = 3903.45
$ws.Cells.Item(129, 14).Value = -13903.45
$ws.Cells.Item(135, 8).Value = 944.93616
$ws.Cells.Item(135, 9).Value = 554.7105
$ws.Cells.Item(135, 10).Value = 2592.5557
$ws.Cells.Item(135, 11).Value = 4992.3945
$ws.Cells.Item(135, 12).Value = 23333.0013
$ws.Cells.Item(135, 13).Value = -2457.3945
$ws.Cells.Item(135, 14).Value = -28403.0013
$ws.Cells.Item(137, 8).Value = 1838.0286
$ws.Cells.Item(137, 9).Value = 1222
$ws.Cells.Item(137, 10).Value = 3182.0908
$ws.Cells.Item(137, 11).Value = 3666
$ws.Cells.Item(137, 12).Value = 9546.2724
$ws.Cells.Item(137, 13).Value = -1116
$ws.Cells.Item(137, 14).Value = -14646.2724
$ws.Cells.Item(138, 8).Value = 1823.6538
$ws.Cells.Item(138, 9).Value = 827.35
$ws.Cells.Item(138, 10).Value = 5144.6665
$ws.Cells.Item(138, 11).Value = 2482.05
$ws.Cells.Item(138, 12).Value = 15433.9995
$ws.Cells.Item(138, 13).Value = 2657.95
$ws.Cells.Item(138, 14).Value = -25713.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(11, 8).Value = 2168
$ws.Cells.Item(11, 9).Value = 500
$ws.Cells.Item(11, 10).Value = 3002
$ws.Cells.Item(11, 11).Value = 500
$ws.Cells.Item(11, 12).Value = 3002
$ws.Cells.Item(11, 13).Value = -356
$ws.Cells.Item(11, 14).Value = -3290
$ws.Cells.Item(132, 8).Value = 2503.205
$ws.Cells.Item(132, 9).Value = 2223.9119
$ws.Cells.Item(132, 10).Value = 4402.4
$ws.Cells.Item(132, 11).Value = 6671.7357
$ws.Cells.Item(132, 12).Value = 13207.2
$ws.Cells.Item(132, 13).Value = -4141.7357
$ws.Cells.Item(132, 14).Value = -18267.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(12, 8).Value = 2824.75
$ws.Cells.Item(12, 9).Value = 499.66666
$ws.Cells.Item(12, 10).Value = 9800
$ws.Cells.Item(12, 11).Value = 499.66666
$ws.Cells.Item(12, 12).Value = 9800
$ws.Cells.Item(12, 13).Value = -331.66666
$ws.Cells.Item(12, 14).Value = -10136
$ws.Cells.Item(24, 8).Value = 4493.2
$ws.Cells.Item(24, 9).Value = 666.5
$ws.Cells.Item(24, 10).Value = 19800
$ws.Cells.Item(24, 11).Value = 666.5
$ws.Cells.Item(24, 12).Value = 19800
$ws.Cells.Item(24, 13).Value = -431.5
$ws.Cells.Item(24, 14).Value = -20270
$ws.Cells.Item(86, 8).Value = 12501689
$ws.Cells.Item(86, 9).Value = 14287359
$ws.Cells.Item(86, 10).Value = 2000
$ws.Cells.Item(86, 11).Value = 14287359
$ws.Cells.Item(86, 12).Value = 2000
$ws.Cells.Item(86, 13).Value = -14286236
$ws.Cells.Item(86, 14).Value = -4246
$ws.Cells.Item(89, 8).Value = 12501689
$ws.Cells.Item(89, 9).Value = 14287359
$ws.Cells.Item(89, 10).Value = 2000
$ws.Cells.Item(89, 11).Value = 71436795
$ws.Cells.Item(89, 12).Value = 10000
$ws.Cells.Item(89, 13).Value = -71431179
$ws.Cells.Item(89, 14).Value = -21232
$ws.Cells.Item(105, 8).Value = 2234.0833
$ws.Cells.Item(105, 9).Value = 1513.625
$ws.Cells.Item(105, 11).Value = 1513.625
$ws.Cells.Item(105, 13).Value = 233.375
$ws.Cells.Item(134, 8).Value = 1925.4546
$ws.Cells.Item(134, 9).Value = 1047.5
$ws.Cells.Item(134, 10).Value = 4266.6665
$ws.Cells.Item(134, 11).Value = 3142.5
$ws.Cells.Item(134, 12).Value = 12799.9995
$ws.Cells.Item(134, 13).Value = -607.5
$ws.Cells.Item(134, 14).Value = -17869.9995

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(12, 8).Value = 1765.3334
$ws.Cells.Item(12, 9).Value = 1765.3334
$ws.Cells.Item(12, 11).Value = 1765.3334
$ws.Cells.Item(12, 13).Value = -1595.3334
$ws.Cells.Item(31, 8).Value = 3866811.2
$ws.Cells.Item(31, 9).Value = 9572575
$ws.Cells.Item(31, 10).Value = 1616.0322
$ws.Cells.Item(31, 11).Value = 9572575
$ws.Cells.Item(31, 12).Value = 1616.0322
$ws.Cells.Item(31, 13).Value = -9572280
$ws.Cells.Item(31, 14).Value = -2206.0322
$ws.Cells.Item(34, 8).Value = 3866811.2
$ws.Cells.Item(34, 9).Value = 9572575
$ws.Cells.Item(34, 10).Value = 1616.0322
$ws.Cells.Item(34, 11).Value = 9572575
$ws.Cells.Item(34, 12).Value = 1616.0322
$ws.Cells.Item(34, 13).Value = -9572373
$ws.Cells.Item(34, 14).Value = -2020.0322

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(19, 8).Value = 820
$ws.Cells.Item(19, 9).Value = 100
$ws.Cells.Item(19, 10).Value = 1000
$ws.Cells.Item(19, 12).Value = 3000
$ws.Cells.Item(19, 13).Value = -126
$ws.Cells.Item(19, 14).Value = -3348
$ws.Cells.Item(100, 8).Value = 10252.934
$ws.Cells.Item(100, 10).Value = 9985.357
$ws.Cells.Item(100, 12).Value = 29956.071
$ws.Cells.Item(100, 14).Value = -31578.071
$ws.Cells.Item(131, 8).Value = 7093065.5
$ws.Cells.Item(131, 9).Value = 253.66667
$ws.Cells.Item(131, 10).Value = 8131038
$ws.Cells.Item(131, 11).Value = 761.00001
$ws.Cells.Item(131, 12).Value = 24393114
$ws.Cells.Item(131, 13).Value = 4278.99999
$ws.Cells.Item(131, 14).Value = -24403194

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(113, 8).Value = 1956.4
$ws.Cells.Item(113, 9).Value = 1840.6666
$ws.Cells.Item(113, 10).Value = 2130
$ws.Cells.Item(113, 11).Value = 1840.6666
$ws.Cells.Item(113, 12).Value = 2130
$ws.Cells.Item(113, 13).Value = 329.3334
$ws.Cells.Item(113, 14).Value = -6470
$ws.Cells.Item(122, 8).Value = 4281.9565
$ws.Cells.Item(122, 9).Value = 3577.8
$ws.Cells.Item(122, 10).Value = 4823.615
$ws.Cells.Item(122, 11).Value = 10733.4
$ws.Cells.Item(122, 12).Value = 14470.845
$ws.Cells.Item(122, 13).Value = -8283.400000000001
$ws.Cells.Item(122, 14).Value = -19370.845

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 1217.6111
$ws.Cells.Item(61, 9).Value = 1315.1428
$ws.Cells.Item(61, 10).Value = 876.25
$ws.Cells.Item(61, 11).Value = 1315.1428
$ws.Cells.Item(61, 12).Value = 876.25
$ws.Cells.Item(61, 13).Value = -1113.1428
$ws.Cells.Item(61, 14).Value = -1280.25
$ws.Cells.Item(68, 8).Value = 56394556
$ws.Cells.Item(68, 9).Value = 112780450
$ws.Cells.Item(68, 10).Value = 8666.666999999999
$ws.Cells.Item(68, 11).Value = 112780450
$ws.Cells.Item(68, 12).Value = 8666.666999999999
$ws.Cells.Item(68, 13).Value = -112779701
$ws.Cells.Item(68, 14).Value = -10164.667
$ws.Cells.Item(71, 8).Value = 56394556
$ws.Cells.Item(71, 9).Value = 112780450
$ws.Cells.Item(71, 10).Value = 8666.666999999999
$ws.Cells.Item(71, 11).Value = 563902250
$ws.Cells.Item(71, 12).Value = 43333.335
$ws.Cells.Item(71, 13).Value = -563898506
$ws.Cells.Item(71, 14).Value = -50821.335
$ws.Cells.Item(113, 8).Value = 1217.6111
$ws.Cells.Item(113, 9).Value = 1315.1428
$ws.Cells.Item(113, 10).Value = 876.25
$ws.Cells.Item(113, 11).Value = 1315.1428
$ws.Cells.Item(113, 12).Value = 876.25
$ws.Cells.Item(113, 13).Value = 854.8571999999999
$ws.Cells.Item(113, 14).Value = -5216.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(11, 8).Value = 20000000
$ws.Cells.Item(11, 9).Value = 20000000
$ws.Cells.Item(11, 11).Value = 20000000
$ws.Cells.Item(11, 13).Value = -19999858
